# Apply the "single choice field" edit to the survey/choices workbook.
$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- survey sheet: add the new "Implementation period" group and the
#     "Monitoring visits" group (single choice field) ---
# (Order of writes matters for the shared-strings table layout: first the
#  "Implementation period" group skeleton, then the "Monitoring visits"
#  group skeleton, then the single-choice question is inserted as a new
#  row inside the first group.)

$survey.Range("A3").Value = "begin group"
$survey.Range("B3").Value = "_1"
$survey.Range("C3").Value = "Implementation period"

$survey.Range("A4").Value = "end group"

$survey.Range("A6").Value = "begin group"
$survey.Range("B6").Value = "_2"
$survey.Range("C6").Value = "Monitoring visits"

$survey.Range("A7").Value = "end group"

# Insert the single-choice question row inside the "Implementation period"
# group, pushing the rest of the rows down by one.
$survey.Rows.Item(4).Insert()

$survey.Range("A4").Value = "select_one IMPLEMENTATION_PERIOD"
$survey.Range("B4").Value = "_1_100"
$survey.Range("C4").Value = "Please select one of the following"

# Column widths to fit new, longer content (matches Excel's computed
# "best fit" widths for the longest entries now present in each column).
$survey.Columns.Item(1).ColumnWidth = 34.8333333333
$survey.Columns.Item(3).ColumnWidth = 30.8333333333

# Page setup (portrait, 300 dpi) for the survey sheet, matching the
# choices sheet that already had these settings.
$survey.PageSetup.Orientation = 1
$survey.PageSetup.HorizontalDpi = 300
$survey.PageSetup.VerticalDpi = 300

# Select A4 and make "survey" the active (visible) sheet/tab.
$survey.Range("A4").Select()
$survey.Activate()

# The "choices" sheet is no longer the active tab.
$choices.Range("B15").Select()

$wb.Save()
